{"js": "// Updates points for hw4:\n//  - \"\u2026\u2026. / 6 Punkten\" (three runs) -> single run \"\u2026\u2026. / 5 Punkten\"\n//  - \"\u2026\u2026. / 4\" -> \"\u2026\u2026. / 5\" (the trailing \" Punkten\" run is left untouched)\n//  - the Word \"last edit\" bookmark (_GoBack) moves from the end of the\n//    document to right after the newly-edited \"\u2026\u2026. / 5\" text\n\n// 1) \"\u2026\u2026. / 6 Punkten\" -> \"\u2026\u2026. / 5 Punkten\" (collapses the 3 runs into 1,\n//    matching what Word itself does when the surviving text carries a\n//    single uniform run -> use insertOoxml so the merged run comes out\n//    identical to a native Word edit, i.e. no stray xml:space marker).\nconst firstMatches = context.document.body.search(\"\u2026\u2026. / 6 Punkten\", { matchCase: true });\nfirstMatches.load(\"text\");\nawait context.sync();\n\nif (firstMatches.items.length > 0) {\n  const target = firstMatches.items[0];\n  const mergedRunOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:rPr>' +\n    '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"IPLCH K+ TU Sans\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n    '<w:bCs/><w:color w:val=\"000000\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/><w:lang w:val=\"de-AT\"/>' +\n    '</w:rPr><w:t>\u2026\u2026. / 5 Punkten</w:t></w:r></w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  target.insertOoxml(mergedRunOoxml, Word.InsertLocation.replace);\n}\n\n// 2) \"\u2026\u2026. / 4\" -> \"\u2026\u2026. / 5\"\nconst secondMatches = context.document.body.search(\"\u2026\u2026. / 4\", { matchCase: true });\nsecondMatches.load(\"text\");\nawait context.sync();\n\nif (secondMatches.items.length > 0) {\n  const target = secondMatches.items[0];\n  target.insertText(\"\u2026\u2026. / 5\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // 3) Move the _GoBack bookmark to right after this edit, like Word does\n  //    when it records the last-edit location.\n  context.document.deleteBookmark(\"_GoBack\");\n  const afterEdit = target.getRange(Word.RangeLocation.end);\n  afterEdit.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Updates points for hw4:\n#  - \"\u2026\u2026. / 6 Punkten\" (three runs) -> single run \"\u2026\u2026. / 5 Punkten\"\n#  - \"\u2026\u2026. / 4\" -> \"\u2026\u2026. / 5\" (the trailing \" Punkten\" run is left untouched)\n#  - the Word \"last edit\" bookmark (_GoBack) moves from the end of the\n#    document to right after the newly-edited \"\u2026\u2026. / 5\" text\n\n$doc = $word.ActiveDocument\n\n# 1) \"\u2026\u2026. / 6 Punkten\" -> \"\u2026\u2026. / 5 Punkten\"\n$rng1 = $doc.Content\nif ($rng1.Find.Execute(\"\u2026\u2026. / 6 Punkten\")) {\n    $rng1.Text = \"\u2026\u2026. / 5 Punkten\"\n}\n\n# 2) \"\u2026\u2026. / 4\" -> \"\u2026\u2026. / 5\", then drop the _GoBack bookmark right after it\n#    (Bookmarks.Add with the existing \"_GoBack\" name relocates it, removing\n#    the old bookmark the same way Word moves its last-edit marker).\n$rng2 = $doc.Content\nif ($rng2.Find.Execute(\"\u2026\u2026. / 4\")) {\n    $rng2.Text = \"\u2026\u2026. / 5\"\n    $rng2.Collapse(0)\n    $doc.Bookmarks.Add(\"_GoBack\", $rng2)\n}\n"}
